# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Femacal de La Calera" / Berenjena
# as a new row 364, pushing the existing rows 364-399 down to 365-400.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 364 (shifts 364..399 -> 365..400)
$ws.Rows.Item(364).Insert()

# Populate the newly inserted row 364 with the new weekly record
$ws.Cells.Item(364, 1).Value = 3
$ws.Cells.Item(364, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(364, 3).Value = "Coquimbo"
$ws.Cells.Item(364, 4).Value = 44946
$ws.Cells.Item(364, 5).Value = 5
$ws.Cells.Item(364, 6).Value = 100112001
$ws.Cells.Item(364, 7).Value = "Berenjena"
$ws.Cells.Item(364, 8).Value = "Sin especificar"
$ws.Cells.Item(364, 9).Value = "Primera"
$ws.Cells.Item(364, 10).Value = 105
$ws.Cells.Item(364, 11).Value = 11000
$ws.Cells.Item(364, 12).Value = 12000
$ws.Cells.Item(364, 13).Value = 11524
$ws.Cells.Item(364, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(364, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(364, 16).Value = 192
$ws.Cells.Item(364, 17).Value = 60
$ws.Cells.Item(364, 18).Value = "Hortaliza"
